# ---- Step 1: insert new sheet '2022-Q3' after '总计' (position 2) ----
$wb = $excel.ActiveWorkbook
$firstSheet = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $firstSheet)
$q3.Name = "2022-Q3"

# helper: true if the string looks like a plain integer/decimal number
function Is-NumericLooking($s) {
    return ($s -match '^[0-9]+(\.[0-9]+)?$')
}

# write a value as TEXT (never auto-converted to a number) unless it is meant to be numeric
function Set-TextCell($cell, $s) {
    if (Is-NumericLooking $s) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $s
}

# ---- Step 2: header row for 2022-Q3 ----
$q3.Cells.Item(1,2).Value = "基金代码"
$q3.Cells.Item(1,3).Value = "基金名称"
$q3.Cells.Item(1,4).Value = "基金规模"
$q3.Cells.Item(1,5).Value = "股票总仓位"
$q3.Cells.Item(1,6).Value = "仓位占比"
$q3.Cells.Item(1,7).Value = "持有市值(亿元)"
$q3.Cells.Item(1,8).Value = "仓位排名"

# ---- Step 3: 36 fund rows for 2022-Q3 (row 2..37) ----
$q3Data = @"
0	720001	财通价值动量混合	38.35	79.13	6.00	2.3010	5
1	001938	中欧时代先锋股票A	128.93	93.47	1.74	2.2434	10
2	003634	嘉实农业产业股票A	22.92	84.79	6.04	1.3844	5
3	001480	财通成长优选混合	20.31	91.20	6.33	1.2856	6
4	005106	银华农业产业股票A	12.00	93.82	5.58	0.6696	6
5	004241	中欧时代先锋股票C	29.58	93.47	1.74	0.5147	10
6	014915	财通匠心优选一年持有期混合A	5.65	81.89	7.69	0.4345	1
7	121005	国投瑞银创新动力混合	11.48	88.74	3.38	0.3880	10
8	015468	嘉实农业产业股票C	6.33	84.79	6.04	0.3823	5
9	013414	太平智远三个月定期开放股票	8.05	88.35	3.99	0.3212	9
10	164403	前海开源沪港深农业混合（LOF）A	3.43	89.98	7.61	0.2610	2
11	001218	国投瑞银精选收益灵活配置混合	3.58	88.17	5.68	0.2033	6
12	121008	国投瑞银成长优选混合	6.02	83.26	3.34	0.2011	9
13	501046	财通多策略福鑫定期开放灵活配置混合	2.82	85.55	6.96	0.1963	4
14	015210	前海开源沪港深农业混合（LOF）C	1.93	89.98	7.61	0.1469	2
15	009062	财通智慧成长混合A	2.17	84.78	6.71	0.1456	4
16	210004	金鹰稳健成长混合	5.18	94.15	2.81	0.1456	7
17	010423	国投瑞银价值成长一年持有期混合A	4.18	87.19	3.35	0.1400	9
18	009063	财通智慧成长混合C	1.50	84.78	6.71	0.1006	4
19	014210	国投瑞银竞争优势混合A	1.48	88.88	3.35	0.0496	9
20	014916	财通匠心优选一年持有期混合C	0.61	81.89	7.69	0.0469	1
21	014064	银华农业产业股票C	0.75	93.82	5.58	0.0418	6
22	002844	金鹰多元策略灵活配置混合	0.48	89.55	8.62	0.0414	3
23	001601	鑫元鑫新收益灵活配置混合A	1.24	62.08	2.98	0.0370	6
24	000845	国投瑞银信息消费灵活配置混合	0.49	88.03	4.99	0.0245	6
25	350009	天治研究驱动混合A	0.29	93.73	8.22	0.0238	6
26	002043	天治研究驱动混合C	0.24	93.73	8.22	0.0197	6
27	001520	国投瑞银研究精选股票	0.53	82.22	3.32	0.0176	9
28	010424	国投瑞银价值成长一年持有期混合C	0.31	87.19	3.35	0.0104	9
29	006522	财通新兴蓝筹混合A	0.22	94.22	3.55	0.0078	10
30	002005	工银新得利混合	0.57	27.58	1.27	0.0072	9
31	011361	华夏博锐一年持有混合（MOM）A	0.11	33.21	2.97	0.0033	4
32	014211	国投瑞银竞争优势混合C	0.05	88.88	3.35	0.0017	9
33	001602	鑫元鑫新收益灵活配置混合C	0.05	62.08	2.98	0.0015	6
34	006523	财通新兴蓝筹混合C	0.03	94.22	3.55	0.0011	10
35	011362	华夏博锐一年持有混合（MOM）C	0.00	33.21	2.97	0	4
"@
$q3Lines = $q3Data -split "`n"
$rowNum = 2
foreach ($line in $q3Lines) {
    $p = $line -split "`t"
    $q3.Cells.Item($rowNum,1).Value = [int]$p[0]
    Set-TextCell $q3.Cells.Item($rowNum,2) $p[1]
    Set-TextCell $q3.Cells.Item($rowNum,3) $p[2]
    Set-TextCell $q3.Cells.Item($rowNum,4) $p[3]
    Set-TextCell $q3.Cells.Item($rowNum,5) $p[4]
    Set-TextCell $q3.Cells.Item($rowNum,6) $p[5]
    if ($p[6] -eq "0") {
        $q3.Cells.Item($rowNum,7).Value = 0
    } else {
        Set-TextCell $q3.Cells.Item($rowNum,7) $p[6]
    }
    $q3.Cells.Item($rowNum,8).Value = [int]$p[7]
    $rowNum = $rowNum + 1
}

# ---- Step 4: update '总计' - shift existing 3 data rows down and insert a 2022-Q3 row on top ----
$total = $wb.Worksheets.Item(1)
for ($r = 4; $r -ge 2; $r--) {
    $src = $r
    $dst = $r + 1
    $total.Range($total.Cells.Item($src,1), $total.Cells.Item($src,4)).Copy($total.Cells.Item($dst,1))
    $total.Cells.Item($dst,1).Value = $dst - 2
}
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q3"
$total.Cells.Item(2,3).Value = 36
$total.Cells.Item(2,4).Value = 11.8
